# Add 2022-Q4 data:
#  - insert a new "2022-Q4" row at the top of the "总计" (summary) sheet
#  - insert a new "2022-Q4" detail worksheet (fund positions) right after "总计"
#  - all pre-existing quarter sheets keep their data, just shift one tab to the right

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) 总计 summary sheet: insert a row under the header and fill in 2022-Q4 totals
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()
$summary.Rows.Item(2).ClearFormats()

$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 1
$summary.Cells.Item(2, 4).Value = 0.05

# re-number the index column (A) sequentially for every data row
for ($r = 2; $r -le 9; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# restore the bordered/centered style on the new A2 cell (Insert+ClearFormats wiped it)
$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(2, 1).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) New "2022-Q4" detail sheet: clone the layout of an existing quarter sheet
#    (so fonts/borders/column layout match exactly) and overwrite its values.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($null, $summary)

# Excel names the clone "<source> (2)" and drops it immediately after the
# sheet passed as "After" above, i.e. right after "总计".
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q4"

$newSheet.Cells.Item(2, 4).Value = "'1.00"
$newSheet.Cells.Item(2, 5).Value = "'78.50"
$newSheet.Cells.Item(2, 6).Value = "'4.86"
$newSheet.Cells.Item(2, 7).Value = "'0.0486"
$newSheet.Cells.Item(2, 8).Value = 6

# the leading apostrophe forces text storage but also tags the cell with a
# "quote prefix" style - strip that back off so the cell format matches the
# plain (unstyled) text cells used throughout the rest of the workbook
$newSheet.Range("D2:G2").Style = "Normal"

# Copy() leaves the freshly cloned sheet focused; restore the tab selection
# to the last sheet ("2021-Q1"), which is where it was originally.
$wb.Worksheets.Item("2021-Q1").Activate()
